$wb = $excel.ActiveWorkbook

# --- Sheet: PIR ---
$ws = $wb.Worksheets.Item('PIR')
$ws.Cells.Item(297, 1).NumberFormat = '@'
$ws.Cells.Item(297, 1).Value = '2026-02-04'
$ws.Cells.Item(297, 2).Value = '14:25:54'
$ws.Cells.Item(297, 3).Value = '14:00'
$ws.Cells.Item(297, 4).Value = 'Bathroom'
$ws.Cells.Item(297, 5).Value = 'No Motion'
$ws.Cells.Item(297, 6).Value = 'Inactive'
$ws.Cells.Item(298, 1).NumberFormat = '@'
$ws.Cells.Item(298, 1).Value = '2026-02-04'
$ws.Cells.Item(298, 2).Value = '14:25:55'
$ws.Cells.Item(298, 3).Value = '14:00'
$ws.Cells.Item(298, 4).Value = 'Bathroom'
$ws.Cells.Item(298, 5).Value = 'Motion Detected'
$ws.Cells.Item(298, 6).Value = 'Active'
$ws.Cells.Item(299, 1).NumberFormat = '@'
$ws.Cells.Item(299, 1).Value = '2026-02-04'
$ws.Cells.Item(299, 2).Value = '14:25:57'
$ws.Cells.Item(299, 3).Value = '14:00'
$ws.Cells.Item(299, 4).Value = 'Bathroom'
$ws.Cells.Item(299, 5).Value = 'No Motion'
$ws.Cells.Item(299, 6).Value = 'Inactive'
$ws.Cells.Item(300, 1).NumberFormat = '@'
$ws.Cells.Item(300, 1).Value = '2026-02-04'
$ws.Cells.Item(300, 2).Value = '14:26:01'
$ws.Cells.Item(300, 3).Value = '14:00'
$ws.Cells.Item(300, 4).Value = 'Bathroom'
$ws.Cells.Item(300, 5).Value = 'No Motion'
$ws.Cells.Item(300, 6).Value = 'Inactive'
$ws.Cells.Item(301, 1).NumberFormat = '@'
$ws.Cells.Item(301, 1).Value = '2026-02-04'
$ws.Cells.Item(301, 2).Value = '14:26:03'
$ws.Cells.Item(301, 3).Value = '14:00'
$ws.Cells.Item(301, 4).Value = 'Bathroom'
$ws.Cells.Item(301, 5).Value = 'Motion Detected'
$ws.Cells.Item(301, 6).Value = 'Active'
$ws.Cells.Item(302, 1).NumberFormat = '@'
$ws.Cells.Item(302, 1).Value = '2026-02-04'
$ws.Cells.Item(302, 2).Value = '14:26:12'
$ws.Cells.Item(302, 3).Value = '14:00'
$ws.Cells.Item(302, 4).Value = 'Bathroom'
$ws.Cells.Item(302, 5).Value = 'No Motion'
$ws.Cells.Item(302, 6).Value = 'Inactive'
$ws.Cells.Item(303, 1).NumberFormat = '@'
$ws.Cells.Item(303, 1).Value = '2026-02-04'
$ws.Cells.Item(303, 2).Value = '14:26:13'
$ws.Cells.Item(303, 3).Value = '14:00'
$ws.Cells.Item(303, 4).Value = 'Bathroom'
$ws.Cells.Item(303, 5).Value = 'Motion Detected'
$ws.Cells.Item(303, 6).Value = 'Active'
$ws.Cells.Item(304, 1).NumberFormat = '@'
$ws.Cells.Item(304, 1).Value = '2026-02-04'
$ws.Cells.Item(304, 2).Value = '14:26:19'
$ws.Cells.Item(304, 3).Value = '14:00'
$ws.Cells.Item(304, 4).Value = 'Bathroom'
$ws.Cells.Item(304, 5).Value = 'No Motion'
$ws.Cells.Item(304, 6).Value = 'Inactive'
$ws.Cells.Item(305, 1).NumberFormat = '@'
$ws.Cells.Item(305, 1).Value = '2026-02-04'
$ws.Cells.Item(305, 2).Value = '14:26:24'
$ws.Cells.Item(305, 3).Value = '14:00'
$ws.Cells.Item(305, 4).Value = 'Bathroom'
$ws.Cells.Item(305, 5).Value = 'No Motion'
$ws.Cells.Item(305, 6).Value = 'Inactive'
$ws.Cells.Item(306, 1).NumberFormat = '@'
$ws.Cells.Item(306, 1).Value = '2026-02-04'
$ws.Cells.Item(306, 2).Value = '14:26:29'
$ws.Cells.Item(306, 3).Value = '14:00'
$ws.Cells.Item(306, 4).Value = 'Bathroom'
$ws.Cells.Item(306, 5).Value = 'No Motion'
$ws.Cells.Item(306, 6).Value = 'Inactive'
$ws.Cells.Item(307, 1).NumberFormat = '@'
$ws.Cells.Item(307, 1).Value = '2026-02-04'
$ws.Cells.Item(307, 2).Value = '14:26:34'
$ws.Cells.Item(307, 3).Value = '14:00'
$ws.Cells.Item(307, 4).Value = 'Bathroom'
$ws.Cells.Item(307, 5).Value = 'No Motion'
$ws.Cells.Item(307, 6).Value = 'Inactive'
$ws.Cells.Item(308, 1).NumberFormat = '@'
$ws.Cells.Item(308, 1).Value = '2026-02-04'
$ws.Cells.Item(308, 2).Value = '14:26:39'
$ws.Cells.Item(308, 3).Value = '14:00'
$ws.Cells.Item(308, 4).Value = 'Bathroom'
$ws.Cells.Item(308, 5).Value = 'No Motion'
$ws.Cells.Item(308, 6).Value = 'Inactive'
$ws.Cells.Item(309, 1).NumberFormat = '@'
$ws.Cells.Item(309, 1).Value = '2026-02-04'
$ws.Cells.Item(309, 2).Value = '14:26:44'
$ws.Cells.Item(309, 3).Value = '14:00'
$ws.Cells.Item(309, 4).Value = 'Bathroom'
$ws.Cells.Item(309, 5).Value = 'No Motion'
$ws.Cells.Item(309, 6).Value = 'Inactive'
$ws.Cells.Item(310, 1).NumberFormat = '@'
$ws.Cells.Item(310, 1).Value = '2026-02-04'
$ws.Cells.Item(310, 2).Value = '14:26:49'
$ws.Cells.Item(310, 3).Value = '14:00'
$ws.Cells.Item(310, 4).Value = 'Bathroom'
$ws.Cells.Item(310, 5).Value = 'No Motion'
$ws.Cells.Item(310, 6).Value = 'Inactive'
$ws.Cells.Item(311, 1).NumberFormat = '@'
$ws.Cells.Item(311, 1).Value = '2026-02-04'
$ws.Cells.Item(311, 2).Value = '14:26:52'
$ws.Cells.Item(311, 3).Value = '14:00'
$ws.Cells.Item(311, 4).Value = 'Bathroom'
$ws.Cells.Item(311, 5).Value = 'Motion Detected'
$ws.Cells.Item(311, 6).Value = 'Active'

# --- Sheet: Humidity ---
$ws = $wb.Worksheets.Item('Humidity')
$ws.Cells.Item(246, 1).NumberFormat = '@'
$ws.Cells.Item(246, 1).Value = '2026-02-04'
$ws.Cells.Item(246, 2).Value = '14:25:55'
$ws.Cells.Item(246, 3).Value = '14:00'
$ws.Cells.Item(246, 4).Value = 'Bathroom'
$ws.Cells.Item(246, 5).NumberFormat = '@'
$ws.Cells.Item(246, 5).Value = '78.8%'
$ws.Cells.Item(246, 6).Value = 'Active'
$ws.Cells.Item(247, 1).NumberFormat = '@'
$ws.Cells.Item(247, 1).Value = '2026-02-04'
$ws.Cells.Item(247, 2).Value = '14:26:05'
$ws.Cells.Item(247, 3).Value = '14:00'
$ws.Cells.Item(247, 4).Value = 'Bathroom'
$ws.Cells.Item(247, 5).NumberFormat = '@'
$ws.Cells.Item(247, 5).Value = '79.8%'
$ws.Cells.Item(247, 6).Value = 'Active'
$ws.Cells.Item(248, 1).NumberFormat = '@'
$ws.Cells.Item(248, 1).Value = '2026-02-04'
$ws.Cells.Item(248, 2).Value = '14:26:10'
$ws.Cells.Item(248, 3).Value = '14:00'
$ws.Cells.Item(248, 4).Value = 'Bathroom'
$ws.Cells.Item(248, 5).NumberFormat = '@'
$ws.Cells.Item(248, 5).Value = '79.7%'
$ws.Cells.Item(248, 6).Value = 'Active'
$ws.Cells.Item(249, 1).NumberFormat = '@'
$ws.Cells.Item(249, 1).Value = '2026-02-04'
$ws.Cells.Item(249, 2).Value = '14:26:15'
$ws.Cells.Item(249, 3).Value = '14:00'
$ws.Cells.Item(249, 4).Value = 'Bathroom'
$ws.Cells.Item(249, 5).NumberFormat = '@'
$ws.Cells.Item(249, 5).Value = '79.7%'
$ws.Cells.Item(249, 6).Value = 'Active'
$ws.Cells.Item(250, 1).NumberFormat = '@'
$ws.Cells.Item(250, 1).Value = '2026-02-04'
$ws.Cells.Item(250, 2).Value = '14:26:30'
$ws.Cells.Item(250, 3).Value = '14:00'
$ws.Cells.Item(250, 4).Value = 'Bathroom'
$ws.Cells.Item(250, 5).NumberFormat = '@'
$ws.Cells.Item(250, 5).Value = '79.6%'
$ws.Cells.Item(250, 6).Value = 'Active'
$ws.Cells.Item(251, 1).NumberFormat = '@'
$ws.Cells.Item(251, 1).Value = '2026-02-04'
$ws.Cells.Item(251, 2).Value = '14:26:35'
$ws.Cells.Item(251, 3).Value = '14:00'
$ws.Cells.Item(251, 4).Value = 'Bathroom'
$ws.Cells.Item(251, 5).NumberFormat = '@'
$ws.Cells.Item(251, 5).Value = '78.7%'
$ws.Cells.Item(251, 6).Value = 'Active'
$ws.Cells.Item(252, 1).NumberFormat = '@'
$ws.Cells.Item(252, 1).Value = '2026-02-04'
$ws.Cells.Item(252, 2).Value = '14:26:41'
$ws.Cells.Item(252, 3).Value = '14:00'
$ws.Cells.Item(252, 4).Value = 'Bathroom'
$ws.Cells.Item(252, 5).NumberFormat = '@'
$ws.Cells.Item(252, 5).Value = '79.5%'
$ws.Cells.Item(252, 6).Value = 'Active'
$ws.Cells.Item(253, 1).NumberFormat = '@'
$ws.Cells.Item(253, 1).Value = '2026-02-04'
$ws.Cells.Item(253, 2).Value = '14:26:46'
$ws.Cells.Item(253, 3).Value = '14:00'
$ws.Cells.Item(253, 4).Value = 'Bathroom'
$ws.Cells.Item(253, 5).NumberFormat = '@'
$ws.Cells.Item(253, 5).Value = '78.6%'
$ws.Cells.Item(253, 6).Value = 'Active'
$ws.Cells.Item(254, 1).NumberFormat = '@'
$ws.Cells.Item(254, 1).Value = '2026-02-04'
$ws.Cells.Item(254, 2).Value = '14:26:51'
$ws.Cells.Item(254, 3).Value = '14:00'
$ws.Cells.Item(254, 4).Value = 'Bathroom'
$ws.Cells.Item(254, 5).NumberFormat = '@'
$ws.Cells.Item(254, 5).Value = '79.5%'
$ws.Cells.Item(254, 6).Value = 'Active'

# --- Sheet: Temperature ---
$ws = $wb.Worksheets.Item('Temperature')
$ws.Cells.Item(246, 1).NumberFormat = '@'
$ws.Cells.Item(246, 1).Value = '2026-02-04'
$ws.Cells.Item(246, 2).Value = '14:25:56'
$ws.Cells.Item(246, 3).Value = '14:00'
$ws.Cells.Item(246, 4).Value = 'Bathroom'
$ws.Cells.Item(246, 5).Value = '24.3C'
$ws.Cells.Item(246, 6).Value = 'Active'
$ws.Cells.Item(247, 1).NumberFormat = '@'
$ws.Cells.Item(247, 1).Value = '2026-02-04'
$ws.Cells.Item(247, 2).Value = '14:26:06'
$ws.Cells.Item(247, 3).Value = '14:00'
$ws.Cells.Item(247, 4).Value = 'Bathroom'
$ws.Cells.Item(247, 5).Value = '24.3C'
$ws.Cells.Item(247, 6).Value = 'Active'
$ws.Cells.Item(248, 1).NumberFormat = '@'
$ws.Cells.Item(248, 1).Value = '2026-02-04'
$ws.Cells.Item(248, 2).Value = '14:26:11'
$ws.Cells.Item(248, 3).Value = '14:00'
$ws.Cells.Item(248, 4).Value = 'Bathroom'
$ws.Cells.Item(248, 5).Value = '24.3C'
$ws.Cells.Item(248, 6).Value = 'Active'
$ws.Cells.Item(249, 1).NumberFormat = '@'
$ws.Cells.Item(249, 1).Value = '2026-02-04'
$ws.Cells.Item(249, 2).Value = '14:26:16'
$ws.Cells.Item(249, 3).Value = '14:00'
$ws.Cells.Item(249, 4).Value = 'Bathroom'
$ws.Cells.Item(249, 5).Value = '24.3C'
$ws.Cells.Item(249, 6).Value = 'Active'
$ws.Cells.Item(250, 1).NumberFormat = '@'
$ws.Cells.Item(250, 1).Value = '2026-02-04'
$ws.Cells.Item(250, 2).Value = '14:26:31'
$ws.Cells.Item(250, 3).Value = '14:00'
$ws.Cells.Item(250, 4).Value = 'Bathroom'
$ws.Cells.Item(250, 5).Value = '24.3C'
$ws.Cells.Item(250, 6).Value = 'Active'
$ws.Cells.Item(251, 1).NumberFormat = '@'
$ws.Cells.Item(251, 1).Value = '2026-02-04'
$ws.Cells.Item(251, 2).Value = '14:26:36'
$ws.Cells.Item(251, 3).Value = '14:00'
$ws.Cells.Item(251, 4).Value = 'Bathroom'
$ws.Cells.Item(251, 5).Value = '24.3C'
$ws.Cells.Item(251, 6).Value = 'Active'
$ws.Cells.Item(252, 1).NumberFormat = '@'
$ws.Cells.Item(252, 1).Value = '2026-02-04'
$ws.Cells.Item(252, 2).Value = '14:26:41'
$ws.Cells.Item(252, 3).Value = '14:00'
$ws.Cells.Item(252, 4).Value = 'Bathroom'
$ws.Cells.Item(252, 5).Value = '24.3C'
$ws.Cells.Item(252, 6).Value = 'Active'
$ws.Cells.Item(253, 1).NumberFormat = '@'
$ws.Cells.Item(253, 1).Value = '2026-02-04'
$ws.Cells.Item(253, 2).Value = '14:26:46'
$ws.Cells.Item(253, 3).Value = '14:00'
$ws.Cells.Item(253, 4).Value = 'Bathroom'
$ws.Cells.Item(253, 5).Value = '24.3C'
$ws.Cells.Item(253, 6).Value = 'Active'
$ws.Cells.Item(254, 1).NumberFormat = '@'
$ws.Cells.Item(254, 1).Value = '2026-02-04'
$ws.Cells.Item(254, 2).Value = '14:26:51'
$ws.Cells.Item(254, 3).Value = '14:00'
$ws.Cells.Item(254, 4).Value = 'Bathroom'
$ws.Cells.Item(254, 5).Value = '24.3C'
$ws.Cells.Item(254, 6).Value = 'Active'
